# Updated symbol list on Sun Jan  1 06:58:47 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for the crypto
# tickers whose market data changed, and re-syncs rows 14-24 / 41-43 whose
# coin listing (Coin name in B, Link in C) shifted by one position versus
# the previous snapshot.
#
# D/E values are plain numeric-looking / percentage-looking text in the
# source data (e.g. "244.06", "-0.41%"), so each literal is written with a
# leading apostrophe to force Excel to store it as text instead of
# auto-converting it to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''244.06'
$ws.Range("E2").Value = '''-0.41%'
# Row 3
$ws.Range("D3").Value = '''26.42'
$ws.Range("E3").Value = '''3.39%'
# Row 4
$ws.Range("D4").Value = '''5.146'
$ws.Range("E4").Value = '''0.63%'
# Row 5
$ws.Range("E5").Value = '''0.40%'
# Row 6
$ws.Range("D6").Value = '''6.469'
$ws.Range("E6").Value = '''-0.06%'
# Row 7
$ws.Range("D7").Value = '''0.8186'
$ws.Range("E7").Value = '''0.03%'
# Row 8
$ws.Range("D8").Value = '''0.8276'
$ws.Range("E8").Value = '''-1.55%'
# Row 9
$ws.Range("E9").Value = '''-0.24%'
# Row 10
$ws.Range("D10").Value = '''0.06923'
$ws.Range("E10").Value = '''-0.40%'
# Row 11
$ws.Range("D11").Value = '''0.02897'
$ws.Range("E11").Value = '''0.78%'
# Row 12
$ws.Range("D12").Value = '''0.09386'
$ws.Range("E12").Value = '''0.06%'
# Row 13
$ws.Range("D13").Value = '''0.001513'
$ws.Range("E13").Value = '''-0.48%'
# Row 14
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = '''0.0005988'
$ws.Range("E14").Value = '''-93.85%'
# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006232'
$ws.Range("E15").Value = '''-0.04%'
# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.655'
$ws.Range("E16").Value = '''3.55%'
# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.033'
$ws.Range("E17").Value = '''0.60%'
# Row 18
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.305'
$ws.Range("E18").Value = '''7.24%'
# Row 19
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3112'
$ws.Range("E19").Value = '''-2.10%'
# Row 20
$ws.Range("B20").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C20").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D20").Value = '''0.03071'
$ws.Range("E20").Value = '''-4.55%'
# Row 21
$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '''0.1299'
$ws.Range("E21").Value = '''-2.25%'
# Row 22
$ws.Range("B22").Value = 'MCDex'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D22").Value = '''3.748'
$ws.Range("E22").Value = '''0.16%'
# Row 23
$ws.Range("B23").Value = 'CoinExToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D23").Value = '''0.04606'
$ws.Range("E23").Value = '''-1.80%'
# Row 24
$ws.Range("B24").Value = 'ZBToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D24").Value = '''0.1342'
$ws.Range("E24").Value = '''-2.45%'
# Row 25
$ws.Range("E25").Value = '''-1.46%'
# Row 26
$ws.Range("D26").Value = '''0.004488'
$ws.Range("E26").Value = '''-2.80%'
# Row 27
$ws.Range("E27").Value = '''-1.08%'
# Row 28
$ws.Range("E28").Value = '''0.62%'
# Row 40
$ws.Range("D40").Value = '''0.03640'
$ws.Range("E40").Value = '''-0.47%'
# Row 41
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '''0.006168'
$ws.Range("E41").Value = '''0.12%'
# Row 42
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1051'
$ws.Range("E42").Value = '''-0.23%'
# Row 43
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002399'
$ws.Range("E43").Value = '''-4.95%'
# Row 44
$ws.Range("E44").Value = '''7.47%'
# Row 45
$ws.Range("D45").Value = '''0.00005364'
$ws.Range("E45").Value = '''0.88%'
# Row 46
$ws.Range("E46").Value = '''-0.05%'
# Row 47
$ws.Range("E47").Value = '''8.21%'
# Row 48
$ws.Range("D48").Value = '''0.002954'
$ws.Range("E48").Value = '''39.15%'
# Row 49
$ws.Range("E49").Value = '''-0.05%'
# Row 50
$ws.Range("E50").Value = '''-0.05%'
